$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: value moved from D2 to C2, with an updated number
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 43.137025930401123

# Row 3: B3 and C3 values removed
$ws.Range("B3:C3").ClearContents()

# Selection now covers B1:E3 instead of B1:AY3
$ws.Range("B1:E3").Select()
